# "Daily 100 Error Counts" weekly upload: append the latest day's numbers
# (11/17/2025) and extend the pre-formatted (but still empty) date column
# down through row 71 so next week's paste keeps the same look.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row for 11/17/2025 (serial 45978)
$ws.Range("A29").Value = 45978
$ws.Range("B29").Value = 1005
$ws.Range("C29").Value = 36
$ws.Range("D29").Value = 969

# Column A uses a date number format. Re-use the format already applied to
# the existing date cells (A2:A27 -> m/d/yyyy) for the last two data rows
# (A28 previously had a different date format) and for the new blank rows
# below the data (A30:A71) so the whole column looks consistent.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A28:A71").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Leave the selection/scroll position where the user was last working
$ws.Range("B30").Select() | Out-Null

Write-Host "done"
